function Set-TextValue {
    param($range, $val)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "60.800.46"

Set-TextValue $ws.Range("D3") "3.373.42"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "569.06"
$ws.Range("E5").Value = "  -1.63%  "

Set-TextValue $ws.Range("D6") "135.91"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("E7").Value = "  +0.09%  "

Set-TextValue $ws.Range("D8") "3.371.10"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("E9").Value = "  -1.35%  "

Set-TextValue $ws.Range("D10") "7.59"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("E11").Value = "  -3.65%  "

Set-TextValue $ws.Range("D12") "0.379"
$ws.Range("E12").Value = "  -2.97%  "

Set-TextValue $ws.Range("D13") "3.944.57"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("E14").Value = "  -0.23%  "

Set-TextValue $ws.Range("D15") "26.04"
$ws.Range("E15").Value = "  +0.18%  "

Set-TextValue $ws.Range("D16") "3.371.88"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("E17").Value = "  -4.37%  "

Set-TextValue $ws.Range("D18") "60.871.65"
$ws.Range("E18").Value = "  -1.43%  "

Set-TextValue $ws.Range("D19") "5.79"
$ws.Range("E19").Value = "  -1.58%  "

Set-TextValue $ws.Range("D20") "13.67"
$ws.Range("E20").Value = "  -4.07%  "

Set-TextValue $ws.Range("D21") "9.21"
$ws.Range("E21").Value = "  -2.58%  "

Set-TextValue $ws.Range("D22") "371.50"
$ws.Range("E22").Value = "  -1.76%  "

Set-TextValue $ws.Range("D23") "3.507.23"
$ws.Range("E23").Value = "  -0.65%  "

Set-TextValue $ws.Range("D24") "0.546"
$ws.Range("E24").Value = "  -2.35%  "

$ws.Range("E25").Value = "  -0.03%  "

Set-TextValue $ws.Range("D26") "70.66"
$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("E27").Value = "  -4.39%  "

Set-TextValue $ws.Range("D28") "0.176"
$ws.Range("E28").Value = "  +9.74%  "

Set-TextValue $ws.Range("D29") "1.59"
$ws.Range("E29").Value = "  -5.57%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.02%  "

Set-TextValue $ws.Range("D31") "7.31"
$ws.Range("E31").Value = "  -3.40%  "

$ws.Range("E32").Value = "  -3.26%  "

$ws.Range("E33").Value = "  -2.91%  "

$ws.Range("E34").Value = "  -0.03%  "

Set-TextValue $ws.Range("D35") "23.25"
$ws.Range("E35").Value = "  -0.90%  "

Set-TextValue $ws.Range("D36") "5.09"
$ws.Range("E36").Value = "  -4.83%  "

Set-TextValue $ws.Range("D37") "1.53"
$ws.Range("E37").Value = "  -2.19%  "

Set-TextValue $ws.Range("D38") "6.75"
$ws.Range("E38").Value = "  -1.70%  "

Set-TextValue $ws.Range("D39") "164.38"
$ws.Range("E39").Value = "  -0.65%  "

Set-TextValue $ws.Range("D40") "0.0755"
$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "1.71"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.771"
$ws.Range("E43").Value = "  -1.56%  "

Set-TextValue $ws.Range("D44") "25.09"
$ws.Range("E44").Value = "  -1.22%  "

Set-TextValue $ws.Range("D45") "41.75"
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("E46").Value = "  -2.62%  "

Set-TextValue $ws.Range("D47") "1.15"
$ws.Range("E47").Value = "  -6.77%  "

Set-TextValue $ws.Range("D48") "2.535.87"
$ws.Range("E48").Value = "  +8.28%  "

Set-TextValue $ws.Range("D49") "23.42"
$ws.Range("E49").Value = "  +2.53%  "

Set-TextValue $ws.Range("D50") "6.74"
$ws.Range("E50").Value = "  -1.85%  "

$ws.Range("E51").Value = "  +0.29%  "
